$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: politeness_score (B34) changes from a text "3" to a true numeric 3
$ws.Range("B34").Value = 3

# New row 35, appended after the old last row (34)
$ws.Range("A35").Value = "Ruilin"

# B35 must stay a text string "3" (not a number) to match the source data
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "3"
$ws.Range("B35").ClearFormats()

$ws.Range("C35").Value = "无"
$ws.Range("D35").Value = "DFT"
$ws.Range("E35").Value = "MET"
$ws.Range("F35").Value = "b3917550-3902-443d-ae6f-4c206bcc883a"
$ws.Range("G35").Value = "HkJ1rgbCb_annotated.xlsx"
$ws.Range("H35").Value = "However, these selections do not seem to directly incorporate the competing/augmenting effects of having different subgraphs within a molecule."
